# Create a test script for a new module named "Opportunity" (VTigerCRM),
# inserted between the existing "LEAD" and "Organization02" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Opportunity" worksheet right after "LEAD" (i.e.
#    right before "Organization02") and make it the active/selected tab.
# ---------------------------------------------------------------------
$leadSheet = $wb.Worksheets.Item("LEAD")
$oppSheet  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $leadSheet)
$oppSheet.Name = "Opportunity"

# ---------------------------------------------------------------------
# 2. Populate the Opportunity sheet with its key/value rows.
# ---------------------------------------------------------------------
$oppSheet.Range("A1").Value = "Opportunity Name"
$oppSheet.Range("B1").Value = "QA Engineer"

$oppSheet.Range("A2").Value = "Related To"
$oppSheet.Range("B2").Value = "Pune01"

$oppSheet.Range("A3").Value = "Type"
$oppSheet.Range("B3").Value = "Existing Business"

$oppSheet.Range("A4").Value = "Lead Source"
$oppSheet.Range("B4").Value = "Existing Customer"

$oppSheet.Range("A5").Value = "Group"
$oppSheet.Range("B5").Value = "Support Group"

$oppSheet.Range("A6").Value = "Sales Stage"
$oppSheet.Range("B6").Value = "Value Proposition"

$oppSheet.Range("A7").Value = "Probability"
$oppSheet.Range("B7").Value = 95

$oppSheet.Range("A8").Value = "Description"
$oppSheet.Range("B8").Value = "Demo purpose"

# B3 ("Existing Business") gets its own distinct (new) font entry.
$oppSheet.Range("B3").Font.Name = "Calibri"
$oppSheet.Range("B3").Font.Size = 11

# ---------------------------------------------------------------------
# 3. Column widths / selection / active cell for the new sheet.
# ---------------------------------------------------------------------
$oppSheet.Columns.Item(1).ColumnWidth = 33.21875
$oppSheet.Columns.Item(2).ColumnWidth = 35.5546875
$oppSheet.Range("B8").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. The LEAD sheet is no longer the selected tab -- Opportunity is.
# ---------------------------------------------------------------------
$oppSheet.Activate()

# ---------------------------------------------------------------------
# 5. Window geometry on the workbook view.
# ---------------------------------------------------------------------
$excel.Width  = 23256
$excel.Height = 13176
$excel.Left   = -108
$excel.Top    = -108

Write-Host "done"
